# Horarios actualizados Linea 141 - 1108
# Refresh the scraped-schedule data across the three worksheets
# (LP1912, LP1912-215, 6203-6173) with the new scrape timestamp
# (04:42:52) and updated "minutes until arrival" countdowns, plus
# newly-scraped rows appended at the bottom of two of the sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 04:42:52"
$ws1.Range("A3").Value = "Total filas: 18"

$ws1.Range("A7").Value = "04:42:52"
$ws1.Range("D7").Value = 3

$ws1.Range("A9").Value = "04:42:52"
$ws1.Range("D9").Value = 11

$ws1.Range("A10").Value = "04:42:52"
$ws1.Range("D10").Value = 34

$ws1.Range("A11").Value = "04:42:52"
$ws1.Range("D11").Value = 39

$ws1.Range("A13").Value = "04:42:52"
$ws1.Range("D13").Value = 52

$ws1.Range("A14").Value = "04:42:52"
$ws1.Range("D14").Value = 64

$ws1.Range("A15").Value = "04:42:52"
$ws1.Range("D15").Value = 71

$ws1.Range("A16").Value = "04:42:52"
$ws1.Range("B16").Value = "06:04"
$ws1.Range("D16").Value = 82

$ws1.Range("B17").Value = "06:05"
$ws1.Range("C17").Value = "16_SANTA ANA"
$ws1.Range("D17").Value = 107

$ws1.Range("A18").Value = "04:42:52"
$ws1.Range("B18").Value = "06:11"
$ws1.Range("C18").Value = "215A_EL PATO"
$ws1.Range("D18").Value = 89

$ws1.Range("A19").Value = "04:42:52"
$ws1.Range("B19").Value = "06:13"
$ws1.Range("C19").Value = "225_HARAS DEL SUR"
$ws1.Range("D19").Value = 91
$ws1.Range("E19").Value = "LP1912"

$ws1.Range("A20").Value = "04:42:52"
$ws1.Range("B20").Value = "06:20"
$ws1.Range("C20").Value = "26_HERNANDEZ"
$ws1.Range("D20").Value = 98
$ws1.Range("E20").Value = "LP1912"

$ws1.Range("A21").Value = "04:42:52"
$ws1.Range("B21").Value = "06:26"
$ws1.Range("C21").Value = "23_HERNANDEZ"
$ws1.Range("D21").Value = 104
$ws1.Range("E21").Value = "LP1912"

$ws1.Range("A22").Value = "04:42:52"
$ws1.Range("B22").Value = "06:29"
$ws1.Range("C22").Value = "86_EST CHICA-ESC AGRARIA"
$ws1.Range("D22").Value = 107
$ws1.Range("E22").Value = "LP1912"

$ws1.Range("A23").Value = "04:42:52"
$ws1.Range("B23").Value = "06:31"
$ws1.Range("C23").Value = "16_SANTA ANA"
$ws1.Range("D23").Value = 109
$ws1.Range("E23").Value = "LP1912"

# ---------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 04:42:52"

$ws2.Range("A6").Value = "04:42:52"
$ws2.Range("D6").Value = 3

$ws2.Range("A8").Value = "04:42:52"
$ws2.Range("D8").Value = 52

$ws2.Range("A9").Value = "04:42:52"
$ws2.Range("D9").Value = 89

# ---------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 04:42:52"
$ws3.Range("A3").Value = "Total filas: 4"

$ws3.Range("A6").Value = "04:42:52"
$ws3.Range("D6").Value = 61

$ws3.Range("A8").Value = "04:42:52"
$ws3.Range("D8").Value = 86

$ws3.Range("A9").Value = "04:42:52"
$ws3.Range("B9").Value = "06:32"
$ws3.Range("C9").Value = "215C_LA PLATA"
$ws3.Range("D9").Value = 110
$ws3.Range("E9").Value = "L6203"
